$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new rows are being inserted (line7/line8) right after the existing
# "line6" row, pushing the "extr1".."extr8" rows down by two. Rows 2-7
# (line1..line6) are untouched.
#
# Create rows 16 and 17 by copying row 15's formatting/structure down first,
# so the new rows end up with identical styling (bold font/border on column A,
# etc.) to the rest of the table.
$ws.Range("A15:E15").Copy($ws.Range("A16:E16"))
$ws.Range("A15:E15").Copy($ws.Range("A17:E17"))

# Target state for rows 8..17 (A,B,C,D,E) after the insertion/shift/fine-tune
$data = @(
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $true),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
